# Time Tracking.xlsx - add Sessions entries for 2024-07-25 through 2024-08-03
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: seed the new unique strings into the shared-string table in the
# exact order they were first typed by the original author (dates are not
# filled top-to-bottom: 07-26 was logged before the 07-25 catch-up entry was
# inserted above it).
# ---------------------------------------------------------------------------
$ws.Range("A1148").Value = "2024-07-26"
$ws.Range("A1146").Value = "2024-07-25"
$ws.Range("A1150").Value = "2024-07-29"
$ws.Range("F1150").Value = "nwragagent v1.0.0"
$ws.Range("A1152").Value = "2024-07-30"
$ws.Range("A1153").Value = "2024-08-01"
$ws.Range("A1155").Value = "2024-08-02"
$ws.Range("A1157").Value = "2024-08-03"

# ---------------------------------------------------------------------------
# Step 2: fill in the rest of row 1146 (2024-07-25, studying AM)
# ---------------------------------------------------------------------------
$ws.Range("B1146").Value = "08:00"
$ws.Range("C1146").Value = "08:45"
$ws.Range("D1146").Value = "0h 45m"
$ws.Range("E1146").Value = "#studying"
$ws.Range("G1146").Value = "'False"
$ws.Range("H1146").Value = "'False"

# row 1147 (2024-07-25, studying PM)
$ws.Range("A1147").Value = "2024-07-25"
$ws.Range("B1147").Value = "16:45"
$ws.Range("C1147").Value = "17:15"
$ws.Range("D1147").Value = "0h 30m"
$ws.Range("E1147").Value = "#studying"
$ws.Range("G1147").Value = "'False"
$ws.Range("H1147").Value = "'False"

# row 1148 (2024-07-26, studying AM)
$ws.Range("B1148").Value = "08:00"
$ws.Range("C1148").Value = "08:45"
$ws.Range("D1148").Value = "0h 45m"
$ws.Range("E1148").Value = "#studying"
$ws.Range("G1148").Value = "'False"
$ws.Range("H1148").Value = "'False"

# row 1149 (2024-07-26, studying PM)
$ws.Range("A1149").Value = "2024-07-26"
$ws.Range("B1149").Value = "16:45"
$ws.Range("C1149").Value = "17:30"
$ws.Range("D1149").Value = "0h 45m"
$ws.Range("E1149").Value = "#studying"
$ws.Range("G1149").Value = "'False"
$ws.Range("H1149").Value = "'False"

# row 1150 (2024-07-29, python / nwragagent release day)
$ws.Range("B1150").Value = "08:00"
$ws.Range("C1150").Value = "15:00"
$ws.Range("D1150").Value = "7h 00m"
$ws.Range("E1150").Value = "#python"
$ws.Range("H1150").Value = "'False"

# row 1151 (2024-07-29, python, evening)
$ws.Range("A1151").Value = "2024-07-29"
$ws.Range("B1151").Value = "20:00"
$ws.Range("C1151").Value = "22:00"
$ws.Range("D1151").Value = "2h 00m"
$ws.Range("E1151").Value = "#python"
$ws.Range("F1151").Value = "nwragagent v1.0.0"
$ws.Range("H1151").Value = "'False"

# row 1152 (2024-07-30, studying)
$ws.Range("B1152").Value = "15:00"
$ws.Range("C1152").Value = "17:30"
$ws.Range("D1152").Value = "2h 30m"
$ws.Range("E1152").Value = "#studying"
$ws.Range("G1152").Value = "'False"
$ws.Range("H1152").Value = "'False"

# row 1153 (2024-08-01, studying AM)
$ws.Range("B1153").Value = "08:00"
$ws.Range("C1153").Value = "08:45"
$ws.Range("D1153").Value = "0h 45m"
$ws.Range("E1153").Value = "#studying"
$ws.Range("G1153").Value = "'False"
$ws.Range("H1153").Value = "'False"

# row 1154 (2024-08-01, studying PM)
$ws.Range("A1154").Value = "2024-08-01"
$ws.Range("B1154").Value = "16:45"
$ws.Range("C1154").Value = "17:30"
$ws.Range("D1154").Value = "0h 45m"
$ws.Range("E1154").Value = "#studying"
$ws.Range("G1154").Value = "'False"
$ws.Range("H1154").Value = "'False"

# row 1155 (2024-08-02, studying AM)
$ws.Range("B1155").Value = "08:15"
$ws.Range("C1155").Value = "08:45"
$ws.Range("D1155").Value = "0h 30m"
$ws.Range("E1155").Value = "#studying"
$ws.Range("G1155").Value = "'False"
$ws.Range("H1155").Value = "'False"

# row 1156 (2024-08-02, studying PM)
$ws.Range("A1156").Value = "2024-08-02"
$ws.Range("B1156").Value = "17:00"
$ws.Range("C1156").Value = "17:45"
$ws.Range("D1156").Value = "0h 45m"
$ws.Range("E1156").Value = "#studying"
$ws.Range("G1156").Value = "'False"
$ws.Range("H1156").Value = "'False"

# row 1157 (2024-08-03, studying afternoon)
$ws.Range("B1157").Value = "16:15"
$ws.Range("C1157").Value = "17:00"
$ws.Range("D1157").Value = "0h 45m"
$ws.Range("E1157").Value = "#studying"
$ws.Range("G1157").Value = "'False"
$ws.Range("H1157").Value = "'False"

# row 1158 (2024-08-03, studying evening)
$ws.Range("A1158").Value = "2024-08-03"
$ws.Range("B1158").Value = "19:30"
$ws.Range("C1158").Value = "20:30"
$ws.Range("D1158").Value = "1h 00m"
$ws.Range("E1158").Value = "#studying"
$ws.Range("G1158").Value = "'False"
$ws.Range("H1158").Value = "'False"

# ---------------------------------------------------------------------------
# Step 3: G1150/G1151 need the literal text "True" WITHOUT the quote-prefix
# formatting that the apostrophe trick above applies everywhere else (mirrors
# a handful of pre-existing cells elsewhere in the sheet, e.g. H732). Copy an
# existing such cell onto them so both value and number format match exactly.
# ---------------------------------------------------------------------------
$ws.Range("H732").Copy($ws.Range("G1150"))
$ws.Range("H732").Copy($ws.Range("G1151"))

# ---------------------------------------------------------------------------
# Step 4: Year/Month helper formulas for the new rows.
# ---------------------------------------------------------------------------
for ($r = 1146; $r -le 1158; $r++) {
    $ws.Range("I$r").Formula = "=YEAR(A$r)"
    $ws.Range("J$r").Formula = "=MONTH(A$r)"
}

# ---------------------------------------------------------------------------
# Step 5: append 12 new blank rows (1160:1171) using the existing blank-row
# formatting as a stamp, matching the author simply extending the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A1159:J1159").Copy($ws.Range("A1160:J1171"))

# ---------------------------------------------------------------------------
# Step 6: restore view state - leave the active cell on the last entry typed.
# ---------------------------------------------------------------------------
$ws.Range("F1157").Select()
